$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row number -> new values for columns B:J (9 values each)
$data = @{
    2 = @(136,112,116,83,150,180,167,62,99)
    3 = @(34,21,19,16,31,39,23,7,12)
    5 = @(0,0,1,0,0,0,0,0,2)
    6 = @(12,12,3,10,13,15,4,4,0)
    7 = @(0,2,0,0,1,1,2,0,0)
    9 = @(10,0,3,2,5,17,2,0,2)
    10 = @(12,7,12,4,12,6,15,3,8)
    11 = @(24,20,23,18,28,24,23,14,24)
    12 = @(4,11,1,11,7,8,11,12,15)
    13 = @(4,3,5,0,0,5,4,0,1)
    14 = @(11,3,12,2,11,4,8,1,5)
    16 = @(4,1,0,5,5,0,0,0,0)
    17 = @(1,2,5,0,5,7,0,1,3)
    19 = @(16,14,17,8,18,29,28,8,14)
    20 = @(2,4,1,0,2,4,4,0,3)
    21 = @(3,2,3,4,5,2,5,2,5)
    22 = @(7,4,10,2,5,17,10,1,3)
    23 = @(2,1,2,1,3,3,6,4,0)
    24 = @(2,3,1,1,3,3,3,1,3)
    25 = @(13,9,7,16,19,20,22,7,7)
    26 = @(1,1,0,1,1,0,0,0,0)
    27 = @(6,1,4,7,8,7,13,1,4)
    29 = @(2,2,0,1,1,1,2,0,0)
    30 = @(4,5,3,7,9,12,7,6,3)
    32 = @(34,37,29,22,42,47,44,2,26)
    35 = @(8,8,3,1,4,4,3,0,1)
    36 = @(6,5,5,3,7,8,9,0,6)
    37 = @(2,2,2,1,0,1,4,0,0)
    38 = @(10,8,8,5,9,17,12,0,3)
    39 = @(8,14,11,12,22,17,16,2,16)
    40 = @(15,11,21,3,12,21,27,24,16)
    43 = @(0,0,3,0,5,8,2,8,5)
    44 = @(9,6,11,3,5,7,13,12,9)
    45 = @(6,1,1,0,0,4,4,4,0)
    46 = @(0,4,6,0,2,2,8,0,0)
}

foreach ($rowKey in $data.Keys) {
    $row = [int]$rowKey
    $values = $data[$rowKey]
    for ($i = 0; $i -lt $values.Count; $i++) {
        $col = $i + 2  # column B is index 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
